$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 7902
$ws.Range("L2").Value = 3431
$ws.Range("L3").Value = 3569
$ws.Range("D4").Value = 1994
$ws.Range("K4").Value = 1769
$ws.Range("L4").Value = 891
$ws.Range("L5").Value = 214
$ws.Range("L6").Value = 3145
$ws.Range("D7").Value = 28185
$ws.Range("K7").Value = 27561
$ws.Range("L7").Value = 11250

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 212
$ws.Range("L3").Value = 239
$ws.Range("L6").Value = 205
$ws.Range("L7").Value = 731

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 87
$ws.Range("L3").Value = 105
$ws.Range("L7").Value = 263

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 147
$ws.Range("L3").Value = 165
$ws.Range("L6").Value = 179
$ws.Range("L7").Value = 526

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 125
$ws.Range("L7").Value = 401

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 76
$ws.Range("L7").Value = 188

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L7").Value = 375
$ws.Range("L8").Value = 731
$ws.Range("L10").Value = 70
$ws.Range("L11").Value = 189
$ws.Range("L13").Value = 17
$ws.Range("L16").Value = 25
$ws.Range("L19").Value = 315
$ws.Range("L20").Value = 287
$ws.Range("L22").Value = 34
$ws.Range("L26").Value = 12
$ws.Range("L27").Value = 104
$ws.Range("L29").Value = 608
$ws.Range("L31").Value = 109
$ws.Range("L33").Value = 526
$ws.Range("L36").Value = 151
$ws.Range("L37").Value = 401
$ws.Range("L40").Value = 29
$ws.Range("L42").Value = 360
$ws.Range("L44").Value = 81
$ws.Range("L48").Value = 159
$ws.Range("L49").Value = 59
$ws.Range("L51").Value = 141
$ws.Range("L52").Value = 235
$ws.Range("L54").Value = 236
$ws.Range("L55").Value = 107
$ws.Range("D63").Value = 374
$ws.Range("K63").Value = 162
$ws.Range("L63").Value = 37
$ws.Range("L64").Value = 75
$ws.Range("L67").Value = 399
$ws.Range("L68").Value = 37
$ws.Range("L69").Value = 30
$ws.Range("L73").Value = 96
$ws.Range("L79").Value = 292
$ws.Range("L80").Value = 35
$ws.Range("L83").Value = 263
$ws.Range("L84").Value = 111
$ws.Range("L85").Value = 565
$ws.Range("L89").Value = 154
$ws.Range("L91").Value = 158
$ws.Range("L94").Value = 135
$ws.Range("L99").Value = 188
$ws.Range("L100").Value = 17
$ws.Range("D101").Value = 28185
$ws.Range("K101").Value = 27561
$ws.Range("L101").Value = 11250

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("L2").Value = 37
$ws.Range("L7").Value = 109

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 117
$ws.Range("L3").Value = 149
$ws.Range("L6").Value = 92
$ws.Range("L7").Value = 399

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 111

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L6").Value = 25
$ws.Range("L7").Value = 59

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L2").Value = 51
$ws.Range("L3").Value = 55
$ws.Range("L6").Value = 111
$ws.Range("L7").Value = 236

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 177
$ws.Range("L3").Value = 238
$ws.Range("L6").Value = 154
$ws.Range("L7").Value = 608

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L3").Value = 37
$ws.Range("L6").Value = 69
$ws.Range("L7").Value = 159

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 107
$ws.Range("L3").Value = 97
$ws.Range("L6").Value = 95
$ws.Range("L7").Value = 315

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L2").Value = 34
$ws.Range("L7").Value = 81

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 107
$ws.Range("L5").Value = 10
$ws.Range("L6").Value = 101
$ws.Range("L7").Value = 360

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("L5").Value = 8
$ws.Range("L6").Value = 17

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L5").Value = 1
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L4").Value = 6
$ws.Range("L7").Value = 107

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 66
$ws.Range("L7").Value = 158

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 100
$ws.Range("L3").Value = 105
$ws.Range("L7").Value = 292

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 75

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L6").Value = 77
$ws.Range("L7").Value = 287

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L6").Value = 41
$ws.Range("L7").Value = 151

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 17

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 121
$ws.Range("L4").Value = 26
$ws.Range("L7").Value = 375

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 32
$ws.Range("L7").Value = 135

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L4").Value = 13
$ws.Range("L7").Value = 189

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L3").Value = 25
$ws.Range("L7").Value = 96

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L2").Value = 33
$ws.Range("L4").Value = 7

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 154

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L6").Value = 29
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L4").Value = 20
$ws.Range("L7").Value = 141

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("L2").Value = 11
$ws.Range("L7").Value = 37

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 162
$ws.Range("L3").Value = 231
$ws.Range("L7").Value = 565

$ws = $wb.Worksheets.Item("Clearing")
$ws.Range("L6").Value = 5
$ws.Range("L7").Value = 34

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 35

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L3").Value = 72
$ws.Range("L6").Value = 61
$ws.Range("L7").Value = 235

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 25
